$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values per cell (case with 380 kV done)
$newValues = @{
    "B2" = 7.080674074869712
    "C2" = 5.959411182113225
    "D2" = 4.475885659471258
    "F2" = 20.88860462921281
    "G2" = 22.86737660241402
    "H2" = 12.87363916951639
    "I2" = 18.61580411706663
    "K2" = 8.507698967479488
    "M2" = 19.97178962817854
    "O2" = 18.7963951308049
    "B3" = 6.773247530042668
    "C3" = 5.900774916668055
    "D3" = 4.383079139917542
    "F3" = 20.92101288002941
    "G3" = 22.9214265074973
    "H3" = 12.9164050755911
    "I3" = 18.70068106776516
    "K3" = 8.337542724873428
    "M3" = 19.37294908440589
    "O3" = 18.8645578292609
    "B4" = 6.575730274252771
    "C4" = 5.864240989951846
    "D4" = 4.324302717674746
    "F4" = 20.94655389911819
    "G4" = 22.96254880407701
    "H4" = 12.94456742240272
    "I4" = 18.75634332163039
    "K4" = 8.230133486267093
    "M4" = 19.00383256119324
    "O4" = 18.9103414451429
    "B5" = 6.493107276536944
    "C5" = 5.849229178715561
    "D5" = 4.299917826627598
    "F5" = 20.95837755525098
    "G5" = 22.98129178862157
    "H5" = 12.95652271272813
    "I5" = 18.77991769310925
    "K5" = 8.185662238633251
    "M5" = 18.85330933974886
    "O5" = 18.92998527395614
    "B6" = 6.479261172016495
    "C6" = 5.846729275642286
    "D6" = 4.29584312076673
    "F6" = 20.96042624811966
    "G6" = 22.9845236413524
    "H6" = 12.95853680889584
    "I6" = 18.78388603347411
    "K6" = 8.178236574564488
    "M6" = 18.82831635868496
    "O6" = 18.9333066490824
    "B7" = 6.574624530045843
    "C7" = 5.864039024529412
    "D7" = 4.323975583549162
    "F7" = 20.94670763068105
    "G7" = 22.96279355418683
    "H7" = 12.94472671637478
    "I7" = 18.75665764465933
    "K7" = 8.229536520929759
    "M7" = 19.00180263486438
    "O7" = 18.91060237650422
    "B8" = 6.976528755862731
    "C8" = 5.939309325908227
    "D8" = 4.444268437423157
    "F8" = 20.89860613991641
    "G8" = 22.88436020588176
    "H8" = 12.88798978964136
    "I8" = 18.64433311580847
    "K8" = 8.449658958193524
    "M8" = 19.76575339036394
    "O8" = 18.81908049082625
    "B9" = 7.692670111394695
    "C9" = 6.082287170524127
    "D9" = 4.665208461883656
    "F9" = 20.8491705600742
    "G9" = 22.79392587563496
    "H9" = 12.79182856930273
    "I9" = 18.452234559021
    "K9" = 8.856553177258377
    "M9" = 21.24210947166765
    "O9" = 18.67088704984858
    "B10" = 8.172172525470796
    "C10" = 6.183970807470414
    "D10" = 4.817517193080981
    "F10" = 20.840348124885
    "G10" = 22.76661142411309
    "H10" = 12.73037444557566
    "I10" = 18.32830313571448
    "K10" = 9.138539953283757
    "M10" = 22.30066478997146
    "O10" = 18.58119539584544
    "B11" = 8.379746930181881
    "C11" = 6.22938810410416
    "D11" = 4.884462273713459
    "F11" = 20.84231988340677
    "G11" = 22.7627643980051
    "H11" = 12.70441222805347
    "I11" = 18.27566652018586
    "K11" = 9.262787417850836
    "M11" = 22.77417491513586
    "O11" = 18.54458268792298
    "B12" = 8.456803925781131
    "C12" = 6.246456724751369
    "D12" = 4.909463206474897
    "F12" = 20.84392716196909
    "G12" = 22.76254608310376
    "H12" = 12.69486759822501
    "I12" = 18.25627307347782
    "K12" = 9.309230363462721
    "M12" = 22.9521510742647
    "O12" = 18.531322671005
    "B13" = 8.44027753587164
    "C13" = 6.242786615038659
    "D13" = 4.904094559682037
    "F13" = 20.84354273792601
    "G13" = 22.76253795501097
    "H13" = 12.69691045354764
    "I13" = 18.26042581469644
    "K13" = 9.299255449884425
    "M13" = 22.91388258965912
    "O13" = 18.53415153785349
    "B14" = 8.386117608115383
    "C14" = 6.230795005812021
    "D14" = 4.886526195720401
    "F14" = 20.84243486994005
    "G14" = 22.76272159267023
    "H14" = 12.70362123957651
    "I14" = 18.27406020328356
    "K14" = 9.266620619374088
    "M14" = 22.78884476430413
    "O14" = 18.54347965288139
    "B15" = 8.352740850360508
    "C15" = 6.223432608439123
    "D15" = 4.875719150896733
    "F15" = 20.84186833298802
    "G15" = 22.76299548566818
    "H15" = 12.70776912727933
    "I15" = 18.2824818716844
    "K15" = 9.246551031762518
    "M15" = 22.71207708803969
    "O15" = 18.54927216824924
    "B16" = 8.158392057510127
    "C16" = 6.180984997948602
    "D16" = 4.81309392654622
    "F16" = 20.84033971518216
    "G16" = 22.76703590374281
    "H16" = 12.73211125528672
    "I16" = 18.33181841599579
    "K16" = 9.130336612227437
    "M16" = 22.26954264874432
    "O16" = 18.58367258793108
    "B17" = 8.036439777745992
    "C17" = 6.154723180545786
    "D17" = 4.774066170699594
    "F17" = 20.84093509469508
    "G17" = 22.77171559382975
    "H17" = 12.74755497931637
    "I17" = 18.3630435153119
    "K17" = 9.057990579757636
    "M17" = 21.99587353234669
    "O17" = 18.60585042712403
    "B18" = 7.965304086186872
    "C18" = 6.139539636321573
    "D18" = 4.751398815761173
    "F18" = 20.84184093785725
    "G18" = 22.7752144570491
    "H18" = 12.75662544333398
    "I18" = 18.38135536726901
    "K18" = 9.016001266447262
    "M18" = 21.83771973996538
    "O18" = 18.61900064906805
    "B19" = 7.941049371969293
    "C19" = 6.134385563070839
    "D19" = 4.743686734324451
    "F19" = 20.84224439396082
    "G19" = 22.77653758485499
    "H19" = 12.75972877369801
    "I19" = 18.38761587146249
    "K19" = 9.001720388869728
    "M19" = 21.78404903687333
    "O19" = 18.62352071237698
    "B20" = 8.049524699673142
    "C20" = 6.157526982130626
    "D20" = 4.778243570069731
    "F20" = 20.84081340411234
    "G20" = 22.77113384827581
    "H20" = 12.74589154809827
    "I20" = 18.35968311214193
    "K20" = 9.0657312430281
    "M20" = 22.02508467908009
    "O20" = 18.60344875329896
    "B21" = 8.402067883267629
    "C21" = 6.234320830764577
    "D21" = 4.891696040586449
    "F21" = 20.84273692462968
    "G21" = 22.76263400921345
    "H21" = 12.7016423378908
    "I21" = 18.27004081888782
    "K21" = 9.276222938921093
    "M21" = 22.82560887383694
    "O21" = 18.54072333890507
    "B22" = 8.623448422486964
    "C22" = 6.283748641798896
    "D22" = 4.963800011510452
    "F22" = 20.8490098569257
    "G22" = 22.7642990161183
    "H22" = 12.6743940365778
    "I22" = 18.21459601506526
    "K22" = 9.410242810935637
    "M22" = 23.34096481118232
    "O22" = 18.50325262977454
    "B23" = 8.506127814706332
    "C23" = 6.257440750767184
    "D23" = 4.925507835303685
    "F23" = 20.84520312871654
    "G23" = 22.76274843051787
    "H23" = 12.68878405016636
    "I23" = 18.2439001394264
    "K23" = 9.339047116765121
    "M23" = 23.06668031722634
    "O23" = 18.52292832061507
    "B24" = 8.043612192793345
    "C24" = 6.156259648343391
    "D24" = 4.776355682517941
    "F24" = 20.84086666509357
    "G24" = 22.77139433769544
    "H24" = 12.74664298837156
    "I24" = 18.36120122877949
    "K24" = 9.062232922432056
    "M24" = 22.01188085320947
    "O24" = 18.60453330463893
    "B25" = 7.506957148250549
    "C25" = 6.044159817287095
    "D25" = 4.607136241390129
    "F25" = 20.85772214081376
    "G25" = 22.81155099385555
    "H25" = 12.8162272312625
    "I25" = 18.50118358880246
    "K25" = 8.74932880491167
    "M25" = 20.84643703892862
    "O25" = 18.70761630282638
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
